# Append run: 2025-09-08 06:26 JST
# Refresh the "ランサーズ" (Lancers) sheet with the latest scrape:
#  - every data row's "取得日時" timestamp advances to the new run time
#  - rows are re-ranked/replaced with the newly scraped listings
#  - the two lowest-priority rows from the previous run fall off the list
#  - column B/H widths are retouched to fit the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Row 2 -----------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(2, 2).Value = "Inkscape「Hershey Text」用svgフォント変換ツール開発(python)"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5389316"
$ws.Cells.Item(2, 7).Value = 315
$ws.Cells.Item(2, 8).Value = "🔥Python ◆ツール,開発"

# --- Row 3 -----------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(3, 2).Value = "Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5273634"
$ws.Cells.Item(3, 7).Value = 298
$ws.Cells.Item(3, 8).Value = "🔥Python ◆開発,スクレイピング"

# --- Row 4 (only the capture timestamp advances) ----------------------
$ws.Cells.Item(4, 1).Value = "2025-09-08 06:26:35"

# --- Row 5 (only the capture timestamp advances) ----------------------
$ws.Cells.Item(5, 1).Value = "2025-09-08 06:26:35"

# --- Row 6 -----------------------------------------------------------
$ws.Cells.Item(6, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(6, 2).Value = "【急募】ECサイトのインタラクティブな商品比較シュミレーションの開発"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5389306"
$ws.Cells.Item(6, 7).Value = 93
$ws.Cells.Item(6, 8).Value = "◆開発 ◇サイト"

# --- Row 7 -----------------------------------------------------------
$ws.Cells.Item(7, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(7, 2).Value = "【急募】Google Cloud WordPress管理画面ログイン設定"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5388922"
$ws.Cells.Item(7, 7).Value = 50
$ws.Cells.Item(7, 8).Value = "◇管理 ○WordPress"

# --- Row 8 (no skill-summary column this time) ------------------------
$ws.Cells.Item(8, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(8, 2).Value = "OR(operations research)にて最適化の仕組みの構築(社内常駐)"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5372984"
$ws.Cells.Item(8, 7).Value = 25
$ws.Range("H8").ClearContents()

# --- Row 9 (no skill-summary column this time) ------------------------
$ws.Cells.Item(9, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(9, 2).Value = "限定公開 PR 限定公開の仕事"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5385681"
$ws.Cells.Item(9, 7).Value = 25
$ws.Range("H9").ClearContents()

# --- Row 10 (no skill-summary column) ----------------------------------
$ws.Cells.Item(10, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(10, 2).Value = "Google Ad Managerの設定支援とGoogleアドセンス・SSPの収益最大化支援"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5389241"
$ws.Cells.Item(10, 7).Value = 13

# --- Row 11 (no skill-summary column) ----------------------------------
$ws.Cells.Item(11, 1).Value = "2025-09-08 06:26:35"
$ws.Cells.Item(11, 2).Value = "【至急】【継続案件】エラーで起動しなくなったエクセルマクロの修正をお願い致します。"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5389081"
$ws.Cells.Item(11, 7).Value = 10

# --- The two oldest/lowest-priority rows drop off the bottom ----------
$ws.Rows("12:13").Delete()

# --- Rebuild the hyperlink collection for the URL column (F) ----------
# (Range.Hyperlinks.Delete clears the whole sheet's collection in this
# host, so drop them all and re-add the eleven that remain.)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5389316")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5273634")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5314730")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5388877")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5389306")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5388922")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5372984")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5385681")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5389241")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5389081")

# --- Column width tweaks (B: 51 -> 50, H: 17 -> 21 OOXML units) -------
$ws.Columns.Item(2).ColumnWidth = 49.17
$ws.Columns.Item(8).ColumnWidth = 20.17
